# Fix MA (moving average) and exponential-smoothing forecast values.
# The workbook stores computed values (no formulas), so we overwrite the
# specific cells whose underlying calculations changed with their corrected
# numeric results, matching the target diff exactly.

$wb = $excel.ActiveWorkbook

# --- moving_average: forecast (D14:D17) now uses the last moving average (C13=99.25) instead of the stale 102 ---
$ws = $wb.Worksheets.Item("moving_average")
$ws.Range("D14").Value2 = 99.25
$ws.Range("D15").Value2 = 99.25
$ws.Range("D16").Value2 = 99.25
$ws.Range("D17").Value2 = 99.25

# --- simple_exponential_smoothing: forecast (D15:D18) now uses the last level (C14=104.364232524118) instead of the stale C13 value ---
$ws = $wb.Worksheets.Item("simple_exponential_smoothing")
$ws.Range("D15").Value2 = 104.364232524118
$ws.Range("D16").Value2 = 104.364232524118
$ws.Range("D17").Value2 = 104.364232524118
$ws.Range("D18").Value2 = 104.364232524118

# --- winter_trendseason: corrected Holt-Winters trend/seasonal smoothing recalculation (C3:M18) ---
$ws = $wb.Worksheets.Item("winter_trendseason")
$ws.Range("C3").Value2 = 118.5388250109192
$ws.Range("D3").Value2 = -2.16754607033665
$ws.Range("E3").Value2 = 0.9521796238547565
$ws.Range("C4").Value2 = 115.9783943403233
$ws.Range("D4").Value2 = -2.206834530362574
$ws.Range("E4").Value2 = 1.032112430639977
$ws.Range("F4").Value2 = 120.521053094399
$ws.Range("G4").Value2 = 4.521053094399022
$ws.Range("H4").Value2 = 4.521053094399022
$ws.Range("I4").Value2 = 27.90090217640276
$ws.Range("J4").Value2 = 5.233819472157904
$ws.Range("K4").Value2 = 3.897459564137088
$ws.Range("L4").Value2 = 4.701778786659649
$ws.Range("C5").Value2 = 114.6402727573129
$ws.Range("D5").Value2 = -2.119963235627355
$ws.Range("E5").Value2 = 0.9633797210446682
$ws.Range("F5").Value2 = 108.7718119528996
$ws.Range("G5").Value2 = -9.22818804710036
$ws.Range("H5").Value2 = 9.22818804710036
$ws.Range("I5").Value2 = 46.98708632848383
$ws.Range("J5").Value2 = 6.56527566380539
$ws.Range("K5").Value2 = 7.820498345000305
$ws.Range("L5").Value2 = 5.741351972773201
$ws.Range("M5").Value2 = 0.1887888583336397
$ws.Range("C6").Value2 = 113.4552765956967
$ws.Range("D6").Value2 = -2.026466528226242
$ws.Range("E6").Value2 = 1.017284241651094
$ws.Range("F6").Value2 = 113.5192491000836
$ws.Range("G6").Value2 = -10.48075089991642
$ws.Range("H6").Value2 = 10.48075089991642
$ws.Range("I6").Value2 = 62.70184960288756
$ws.Range("J6").Value2 = 7.544144472833146
$ws.Range("K6").Value2 = 8.452218467674529
$ws.Range("L6").Value2 = 6.419068596498533
$ws.Range("M6").Value2 = -1.224963285894029
$ws.Range("C7").Value2 = 110.474135341901
$ws.Range("D7").Value2 = -2.121934000783193
$ws.Range("E7").Value2 = 0.9438598297173457
$ws.Range("F7").Value2 = 106.1002424566271
$ws.Range("G7").Value2 = 10.10024245662714
$ws.Range("H7").Value2 = 10.10024245662714
$ws.Range("I7").Value2 = 70.56445921884075
$ws.Range("J7").Value2 = 8.055364069591946
$ws.Range("K7").Value2 = 10.52108589231994
$ws.Range("L7").Value2 = 7.239472055662816
$ws.Range("M7").Value2 = 0.1066298737717619
$ws.Range("C8").Value2 = 108.9772797052805
$ws.Range("D8").Value2 = -2.059426164366925
$ws.Range("E8").Value2 = 1.038098260875868
$ws.Range("F8").Value2 = 111.8316538913732
$ws.Range("G8").Value2 = -7.168346108626807
$ws.Range("H8").Value2 = 7.168346108626807
$ws.Range("I8").Value2 = 67.36791367121147
$ws.Range("J8").Value2 = 7.907527742764422
$ws.Range("K8").Value2 = 6.023820259350257
$ws.Range("L8").Value2 = 7.036863422944056
$ws.Range("M8").Value2 = -0.7978983899833784
$ws.Range("C9").Value2 = 106.2636729940974
$ws.Range("D9").Value2 = -2.12484421904854
$ws.Range("E9").Value2 = 0.9573830643635187
$ws.Range("F9").Value2 = 103.00249191894
$ws.Range("G9").Value2 = 7.002491918939967
$ws.Range("H9").Value2 = 7.002491918939967
$ws.Range("I9").Value2 = 64.74891072886977
$ws.Range("J9").Value2 = 7.7782369107895
$ws.Range("K9").Value2 = 7.294262415562465
$ws.Range("L9").Value2 = 7.073634707603829
$ws.Range("M9").Value2 = 0.08910608820334663
$ws.Range("C10").Value2 = 103.7903607298328
$ws.Range("D10").Value2 = -2.159691023570148
$ws.Range("E10").Value2 = 1.01383083963906
$ws.Range("F10").Value2 = 105.9387894568587
$ws.Range("G10").Value2 = 3.938789456858686
$ws.Range("H10").Value2 = 3.938789456858686
$ws.Range("I10").Value2 = 58.59455468594369
$ws.Range("J10").Value2 = 7.298305979048148
$ws.Range("K10").Value2 = 3.861558291037928
$ws.Range("L10").Value2 = 6.67212515553309
$ws.Range("M10").Value2 = 0.6346510730565609
$ws.Range("C11").Value2 = 103.1634612926596
$ws.Range("D11").Value2 = -2.006411864930449
$ws.Range("E11").Value2 = 0.9580394169548856
$ws.Range("F11").Value2 = 95.92510660301284
$ws.Range("G11").Value2 = -16.07489339698716
$ws.Range("H11").Value2 = 16.07489339698716
$ws.Range("I11").Value2 = 80.79540391245013
$ws.Range("J11").Value2 = 8.27348235881915
$ws.Range("K11").Value2 = 14.35258339016711
$ws.Range("L11").Value2 = 7.525509403825759
$ws.Range("M11").Value2 = -1.383095434257066
$ws.Range("C12").Value2 = 100.8960084955011
$ws.Range("D12").Value2 = -2.032515958153256
$ws.Range("E12").Value2 = 1.035382622279899
$ws.Range("F12").Value2 = 105.0109570862599
$ws.Range("G12").Value2 = 3.010957086259864
$ws.Range("H12").Value2 = 3.010957086259864
$ws.Range("I12").Value2 = 73.62244977873496
$ws.Range("J12").Value2 = 7.747229831563222
$ws.Range("K12").Value2 = 2.951918712019475
$ws.Range("L12").Value2 = 7.06815033464513
$ws.Range("M12").Value2 = -1.088396597616864
$ws.Range("C13").Value2 = 98.61435400712293
$ws.Range("D13").Value2 = -2.057429811175746
$ws.Range("E13").Value2 = 0.954937464578528
$ws.Range("F13").Value2 = 94.65023343908592
$ws.Range("G13").Value2 = 2.650233439085923
$ws.Range("H13").Value2 = 2.650233439085923
$ws.Range("I13").Value2 = 67.56802136990899
$ws.Range("J13").Value2 = 7.283866523156195
$ws.Range("K13").Value2 = 2.880688520745569
$ws.Range("L13").Value2 = 6.687471987926989
$ws.Range("M13").Value2 = -0.7937851595937833
$ws.Range("C14").Value2 = 95.94507175124384
$ws.Range("D14").Value2 = -2.11861505564608
$ws.Range("E14").Value2 = 1.00729369026983
$ws.Range("F14").Value2 = 97.89238753054219
$ws.Range("G14").Value2 = 6.892387530542194
$ws.Range("H14").Value2 = 6.892387530542194
$ws.Range("I14").Value2 = 65.89610341168104
$ws.Range("J14").Value2 = 7.251243273771695
$ws.Range("K14").Value2 = 7.574052231365049
$ws.Range("L14").Value2 = 6.761353674880159
$ws.Range("M14").Value2 = 0.1531547540289855
$ws.Range("E15").Value2 = 0.9580394169548856
$ws.Range("F15").Value2 = 94.47643821636673
$ws.Range("E16").Value2 = 1.035382622279899
$ws.Range("F16").Value2 = 99.97336137328116
$ws.Range("E17").Value2 = 0.954937464578528
$ws.Range("F17").Value2 = 90.24110757174648
$ws.Range("E18").Value2 = 1.00729369026983
$ws.Range("F18").Value2 = 93.11630836049916
